$wb = $excel.ActiveWorkbook

# --- ALC!row6 (hunk 0) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 250241.25

# --- ALC!row100 (hunk 1) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 5130.643
$ws.Range("I100").Value = 6218.8423
$ws.Range("K100").Value = 6218.8423
$ws.Range("M100").Value = -5677.8423

# --- ALC!row137 (hunk 2) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1015.875
$ws.Range("I137").Value = 901.64105
$ws.Range("J137").Value = 1510.8889
$ws.Range("K137").Value = 2704.92315
$ws.Range("L137").Value = 4532.6667
$ws.Range("M137").Value = -154.9231499999996
$ws.Range("N137").Value = -9632.6667

# --- ALC!row138 (hunk 3) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2269.0576
$ws.Range("I138").Value = 1424.8695
$ws.Range("J138").Value = 2938.5862
$ws.Range("K138").Value = 4274.6085
$ws.Range("L138").Value = 8815.758600000001
$ws.Range("M138").Value = 865.3914999999997
$ws.Range("N138").Value = -19095.7586

# --- ARM!row6 (hunk 4) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 2000
$ws.Range("I6").Value = 2000
$ws.Range("K6").Value = 2000
$ws.Range("M6").Value = -1827

# --- ARM!row61 (hunk 5) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1366.8889
$ws.Range("I61").Value = 1108.9565
$ws.Range("J61").Value = 2850
$ws.Range("K61").Value = 1108.9565
$ws.Range("L61").Value = 2850
$ws.Range("M61").Value = -896.9565
$ws.Range("N61").Value = -3274

# --- ARM!row109 (hunk 6) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H109").Value = 31310.75
$ws.Range("J109").Value = 31310.75
$ws.Range("L109").Value = 31310.75
$ws.Range("N109").Value = -34084.75

# --- ARM!row136 (hunk 7) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 1366.8889
$ws.Range("I136").Value = 1108.9565
$ws.Range("J136").Value = 2850
$ws.Range("K136").Value = 3326.8695
$ws.Range("L136").Value = 8550
$ws.Range("M136").Value = -776.8694999999998
$ws.Range("N136").Value = -13650

# --- BSM!row134 (hunk 8) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3062.0981
$ws.Range("I134").Value = 794.6047
$ws.Range("J134").Value = 15249.875
$ws.Range("K134").Value = 2383.8141
$ws.Range("L134").Value = 45749.625
$ws.Range("M134").Value = 151.1858999999999
$ws.Range("N134").Value = -50819.625

# --- CRP!row3 (hunk 9) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 1500
$ws.Range("I3").Value = 1500
$ws.Range("K3").Value = 1500
$ws.Range("M3").Value = -1387

# --- CRP!row31 (hunk 10) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1440.4375
$ws.Range("I31").Value = 1138.4
$ws.Range("J31").Value = 1577.7273
$ws.Range("K31").Value = 1138.4
$ws.Range("L31").Value = 1577.7273
$ws.Range("M31").Value = -843.4000000000001
$ws.Range("N31").Value = -2167.7273

# --- CRP!row34 (hunk 11) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 1440.4375
$ws.Range("I34").Value = 1138.4
$ws.Range("J34").Value = 1577.7273
$ws.Range("K34").Value = 1138.4
$ws.Range("L34").Value = 1577.7273
$ws.Range("M34").Value = -936.4000000000001
$ws.Range("N34").Value = -1981.7273

# --- CRP!row58 (hunk 12) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 13158657
$ws.Range("I58").Value = 17242104
$ws.Range("J58").Value = 880.1111
$ws.Range("K58").Value = 17242104
$ws.Range("L58").Value = 880.1111
$ws.Range("M58").Value = -17241901
$ws.Range("N58").Value = -1286.1111

# --- CRP!row132 (hunk 13) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 6174039.5
$ws.Range("I132").Value = 865.37836
$ws.Range("J132").Value = 19609772
$ws.Range("K132").Value = 2596.13508
$ws.Range("L132").Value = 58829316
$ws.Range("M132").Value = -66.13508000000002
$ws.Range("N132").Value = -58834376

# --- CRP!row134 (hunk 14) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 12821304
$ws.Range("I134").Value = 765.0303
$ws.Range("K134").Value = 2295.0909
$ws.Range("M134").Value = 239.9090999999999

# --- CRP!row136 (hunk 15) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 13158657
$ws.Range("I136").Value = 17242104
$ws.Range("J136").Value = 880.1111
$ws.Range("K136").Value = 51726312
$ws.Range("L136").Value = 2640.3333
$ws.Range("M136").Value = -51723762
$ws.Range("N136").Value = -7740.3333

# --- CRP!row141 (hunk 16) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H141").Value = 33738.11
$ws.Range("J141").Value = 33738.11
$ws.Range("L141").Value = 33738.11
$ws.Range("N141").Value = -44098.11

# --- CUL!row4 (hunk 17) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 200268.2
$ws.Range("I4").Value = 250312.75
$ws.Range("K4").Value = 750938.25
$ws.Range("M4").Value = -750826.25

# --- CUL!row131 (hunk 18) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 817.03
$ws.Range("I131").Value = 488.16666
$ws.Range("J131").Value = 838.0213
$ws.Range("K131").Value = 1464.49998
$ws.Range("L131").Value = 2514.0639
$ws.Range("M131").Value = 3575.50002
$ws.Range("N131").Value = -12594.0639

# --- CUL!row138 (hunk 19) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H138").Value = 17544820
$ws.Range("I138").Value = 17544820
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 52634460
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = -52629320
$ws.Range("M138").ClearContents()

# --- CUL!row141 (hunk 20) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H141").Value = 3236.7856
$ws.Range("I141").Value = 2641.8
$ws.Range("J141").Value = 4724.25
$ws.Range("K141").Value = 7925.400000000001
$ws.Range("L141").Value = 14172.75
$ws.Range("M141").Value = -2745.400000000001
$ws.Range("N141").Value = -24532.75

# --- GSM!row5 (hunk 21) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("I5").Value = 1000
$ws.Range("K5").Value = 1000
$ws.Range("M5").Value = -888

# --- LTW!row40 (hunk 22) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7354705
$ws.Range("I40").Value = 1389.2693
$ws.Range("J40").Value = 31252982
$ws.Range("K40").Value = 1389.2693
$ws.Range("L40").Value = 31252982
$ws.Range("M40").Value = -1253.2693
$ws.Range("N40").Value = -31253254

# --- LTW!row43 (hunk 23) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H43").Value = 4000
$ws.Range("J43").Value = 4000
$ws.Range("L43").Value = 4000
$ws.Range("N43").Value = -4386

# --- LTW!row132 (hunk 24) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 18187516
$ws.Range("I132").Value = 38463292
$ws.Range("J132").Value = 9235.585999999999
$ws.Range("K132").Value = 115389876
$ws.Range("L132").Value = 27706.758
$ws.Range("M132").Value = -115387346
$ws.Range("N132").Value = -32766.758

# --- WVR!row2 (hunk 25) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50000
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()

# --- WVR!row5 (hunk 26) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 35626
$ws.Range("J5").Value = 35626
$ws.Range("L5").Value = 35626
$ws.Range("N5").Value = -35850

# --- WVR!row37 (hunk 27) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 6000
$ws.Range("J37").Value = 6000
$ws.Range("L37").Value = 6000
$ws.Range("N37").Value = -6406

# --- WVR!row132 (hunk 28) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 18167.836
$ws.Range("I132").Value = 24897.861
$ws.Range("J132").Value = 6109.875
$ws.Range("K132").Value = 74693.583
$ws.Range("L132").Value = 18329.625
$ws.Range("M132").Value = -72163.583
$ws.Range("N132").Value = -23389.625
